$d = $word.ActiveDocument

# Locate the paragraph that holds the pricing figure "$40,000" - the new
# "Previous Experience with Wannon Water" section is inserted right after it
# (and before the existing blank paragraph / "Please send any questions to:"
# block).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("`$40,000", $true, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '`$40,000' paragraph to anchor the new section."
}
$pricingPara = $searchRange.Paragraphs(1)

$insertionPoint = $d.Range($pricingPara.Range.End, $pricingPara.Range.End)

$newBodyXml = @'
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Previous Experience with Wannon Water</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Companny</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">XX is looking forward to establishing a successful relationship with </w:t>
  </w:r>
  <w:r>
    <w:t>Wannon Water.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">  We have worked closely with numerous other Victorian Water businesses.</w:t>
  </w:r>
</w:p>
<w:p/>
'@

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newBodyXml + '</w:body></w:wordDocument>'

$insertionPoint.InsertXML($packageXml) | Out-Null
